# "Refined metadata to be additional tab"
#
# 1. Refreshes the `time_taken` (column F) timestamps on the "data" sheet
#    to reflect a later export run.
# 2. Adds a new "metadata" worksheet (after "data") describing the panel
#    export itself (name/id/version/version-created/query-time/request url).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1. Update the F column ("time_taken") timestamps on "data" ---------
$ws.Range("F2").Value = "2021-10-05 14:19:59.968195"
$ws.Range("F3").Value = "2021-10-05 14:19:59.968203"
$ws.Range("F4").Value = "2021-10-05 14:19:59.968206"
$ws.Range("F5").Value = "2021-10-05 14:19:59.968209"
$ws.Range("F6").Value = "2021-10-05 14:19:59.968212"
$ws.Range("F7").Value = "2021-10-05 14:19:59.968214"
$ws.Range("F8").Value = "2021-10-05 14:19:59.968217"
$ws.Range("F9").Value = "2021-10-05 14:19:59.968219"
$ws.Range("F10").Value = "2021-10-05 14:19:59.968222"
$ws.Range("F11").Value = "2021-10-05 14:19:59.968225"
$ws.Range("F12").Value = "2021-10-05 14:19:59.968227"
$ws.Range("F13").Value = "2021-10-05 14:19:59.968230"
$ws.Range("F14").Value = "2021-10-05 14:19:59.968233"
$ws.Range("F15").Value = "2021-10-05 14:19:59.968235"
$ws.Range("F16").Value = "2021-10-05 14:19:59.968238"
$ws.Range("F17").Value = "2021-10-05 14:19:59.968240"
$ws.Range("F18").Value = "2021-10-05 14:19:59.968243"
$ws.Range("F19").Value = "2021-10-05 14:19:59.968247"
$ws.Range("F20").Value = "2021-10-05 14:19:59.968249"
$ws.Range("F21").Value = "2021-10-05 14:19:59.968252"
$ws.Range("F22").Value = "2021-10-05 14:19:59.968255"
$ws.Range("F23").Value = "2021-10-05 14:19:59.968258"
$ws.Range("F24").Value = "2021-10-05 14:19:59.968260"
$ws.Range("F25").Value = "2021-10-05 14:19:59.968263"
$ws.Range("F26").Value = "2021-10-05 14:19:59.968266"
$ws.Range("F27").Value = "2021-10-05 14:19:59.968269"
$ws.Range("F28").Value = "2021-10-05 14:19:59.968272"
$ws.Range("F29").Value = "2021-10-05 14:19:59.968274"
$ws.Range("F30").Value = "2021-10-05 14:19:59.968277"
$ws.Range("F31").Value = "2021-10-05 14:19:59.968279"
$ws.Range("F32").Value = "2021-10-05 14:19:59.968282"
$ws.Range("F33").Value = "2021-10-05 14:19:59.968284"
$ws.Range("F34").Value = "2021-10-05 14:19:59.968287"
$ws.Range("F35").Value = "2021-10-05 14:19:59.968290"
$ws.Range("F36").Value = "2021-10-05 14:19:59.968293"
$ws.Range("F37").Value = "2021-10-05 14:19:59.968295"
$ws.Range("F38").Value = "2021-10-05 14:19:59.968298"
$ws.Range("F39").Value = "2021-10-05 14:19:59.968300"
$ws.Range("F40").Value = "2021-10-05 14:19:59.968303"
$ws.Range("F41").Value = "2021-10-05 14:19:59.968306"
$ws.Range("F42").Value = "2021-10-05 14:19:59.968309"
$ws.Range("F43").Value = "2021-10-05 14:19:59.968311"
$ws.Range("F44").Value = "2021-10-05 14:19:59.968314"
$ws.Range("F45").Value = "2021-10-05 14:19:59.968317"
$ws.Range("F46").Value = "2021-10-05 14:19:59.968319"
$ws.Range("F47").Value = "2021-10-05 14:19:59.968322"
$ws.Range("F48").Value = "2021-10-05 14:19:59.968325"
$ws.Range("F49").Value = "2021-10-05 14:19:59.968327"
$ws.Range("F50").Value = "2021-10-05 14:19:59.968330"
$ws.Range("F51").Value = "2021-10-05 14:19:59.968332"
$ws.Range("F52").Value = "2021-10-05 14:19:59.968335"
$ws.Range("F53").Value = "2021-10-05 14:19:59.968338"
$ws.Range("F54").Value = "2021-10-05 14:19:59.968341"
$ws.Range("F55").Value = "2021-10-05 14:19:59.968344"
$ws.Range("F56").Value = "2021-10-05 14:19:59.968346"
$ws.Range("F57").Value = "2021-10-05 14:19:59.968349"
$ws.Range("F58").Value = "2021-10-05 14:19:59.968351"
$ws.Range("F59").Value = "2021-10-05 14:19:59.968354"

# --- 2. Add the new "metadata" sheet, positioned right after "data" -----
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1) - reuse the bold/bordered header style from "data"!B1
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Row index cell (A2) - reuse the style from "data"!A2
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Dilated cardiomyopathy - adult and teen"
$meta.Range("C2").Value = 652

# data_version must stay a text value ("1.25"), not become numeric 1.25.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.25"
$meta.Range("D2").ClearFormats()

$meta.Range("E2").Value = "2021-05-12T14:06:03.417132Z"
$meta.Range("F2").Value = "2021-10-05 14:19:59.964904"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/652/?format=json"

# Keep "data" as the active/selected sheet, matching the original workbook.
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
